$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.278.39'
$ws.Range('E2').Value = '  -0.07%  '

$ws.Range('D3').Value = '1.685.81'
$ws.Range('E3').Value = '  +0.45%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.15%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5320'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.19%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.008'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2712'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.90%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06406'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.03%  '

$ws.Range('E10').Value = '  -0.79%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07675'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.35%  '

$ws.Range('D12').Value = '1.688.18'
$ws.Range('E12').Value = '  +0.42%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.528'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.34%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5786'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.27%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008363'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.36%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.86'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.43%  '

$ws.Range('D17').Value = '26.309.21'
$ws.Range('E17').Value = '  -0.05%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.008'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.06%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.892'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.52%  '

$ws.Range('E20').Value = '  -0.17%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.73'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.11%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.253'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.13%  '

$ws.Range('E23').Value = '  -0.07%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.18'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.97%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.829'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.59%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1270'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.32%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.81'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.39%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06255'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.92%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.373'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.61%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.323'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.08%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.598'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.49%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.576'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.15%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.690'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.19%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.030'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.32%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6176'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.09%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.429'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.80%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.754'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.52%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.240'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.66%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01636'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.87%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8986'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.25%  '

$ws.Range('D41').Value = '1.108.42'
$ws.Range('E41').Value = '  -0.69%  '

$ws.Range('E42').Value = '  -0.31%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.96'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.46%  '

$ws.Range('D44').Value = '1.838.52'
$ws.Range('E44').Value = '  +0.64%  '

$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000114'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.68%  '

$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.66'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.34%  '

$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.008'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.03%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.087'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.13%  '

$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05280'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.28%  '

$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4293'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.09%  '

$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.053'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.04%  '
